$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "272.35"
$ws.Range("D2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.09"
$ws.Range("D3").NumberFormat = "General"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.253"
$ws.Range("D4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06197"
$ws.Range("D5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.561"
$ws.Range("D6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.543"
$ws.Range("D7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.431"
$ws.Range("D8").NumberFormat = "General"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8231"
$ws.Range("D9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1649"
$ws.Range("D10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08285"
$ws.Range("D11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03547"
$ws.Range("D12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03184"
$ws.Range("D13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09181"
$ws.Range("D14").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001628"
$ws.Range("D16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04668"
$ws.Range("D17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006299"
$ws.Range("D18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006194"
$ws.Range("D19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001067"
$ws.Range("D20").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.725"
$ws.Range("D22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.226"
$ws.Range("D23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01380"
$ws.Range("D24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3291"
$ws.Range("D25").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04704"
$ws.Range("D40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007026"
$ws.Range("D41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004500"
$ws.Range("D42").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01050"
$ws.Range("D44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006376"
$ws.Range("D45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009900"
$ws.Range("D46").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9905"
$ws.Range("D48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001395"
$ws.Range("D49").NumberFormat = "General"
